$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 01:04"

# Paraguay/Kirguistan swap order (rows 67-68) as data refresh changed relative ranking
$ws.Range("A67").Value = "Paraguay"
$ws.Range("A68").Value = "Kirguistan"

# Uruguay/Burkina Faso swap order (rows 157-158) as data refresh changed relative ranking
$ws.Range("A157").Value = "Uruguay"
$ws.Range("A158").Value = "Burkina Faso"

# Update refreshed COVID statistics for affected countries (columns B:H)
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 7942017
$ws.Cells.Item(4, 3).Value = 47388
$ws.Cells.Item(4, 4).Value = 5084615
$ws.Cells.Item(4, 5).Value = 2638158
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 597
$ws.Cells.Item(4, 8).Value = 219244

# Row 6: Brasil
$ws.Cells.Item(6, 2).Value = 5082637
$ws.Cells.Item(6, 3).Value = 25447
$ws.Cells.Item(6, 4).Value = 4453722
$ws.Cells.Item(6, 5).Value = 478717
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 506
$ws.Cells.Item(6, 8).Value = 150198

# Row 8: Colombia
$ws.Cells.Item(8, 2).Value = 902747
$ws.Cells.Item(8, 3).Value = 8447
$ws.Cells.Item(8, 4).Value = 783131
$ws.Cells.Item(8, 5).Value = 91956
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 165
$ws.Cells.Item(8, 8).Value = 27660

# Row 10: Argentina
$ws.Cells.Item(10, 2).Value = 883882
$ws.Cells.Item(10, 3).Value = 12414
$ws.Cells.Item(10, 4).Value = 709464
$ws.Cells.Item(10, 5).Value = 150837
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 356
$ws.Cells.Item(10, 8).Value = 23581

# Row 45: Egipto
$ws.Cells.Item(45, 2).Value = 104387
$ws.Cells.Item(45, 3).Value = 125
$ws.Cells.Item(45, 4).Value = 97643
$ws.Cells.Item(45, 5).Value = 704
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 11
$ws.Cells.Item(45, 8).Value = 6040

# Row 57: Barein
$ws.Cells.Item(57, 2).Value = 75287
$ws.Cells.Item(57, 3).Value = 427
$ws.Cells.Item(57, 4).Value = 70808
$ws.Cells.Item(57, 5).Value = 4206
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 2
$ws.Cells.Item(57, 8).Value = 273

# Row 61: Nigeria
$ws.Cells.Item(61, 2).Value = 60103
$ws.Cells.Item(61, 3).Value = 111
$ws.Cells.Item(61, 4).Value = 51711
$ws.Cells.Item(61, 5).Value = 7277
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 1115

# Row 65: Argelia
$ws.Cells.Item(65, 2).Value = 52940
$ws.Cells.Item(65, 3).Value = 136
$ws.Cells.Item(65, 4).Value = 37170
$ws.Cells.Item(65, 5).Value = 13975
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 6
$ws.Cells.Item(65, 8).Value = 1795

# Row 67: Paraguay
$ws.Cells.Item(67, 2).Value = 48978
$ws.Cells.Item(67, 3).Value = 703
$ws.Cells.Item(67, 4).Value = 31351
$ws.Cells.Item(67, 5).Value = 16562
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 20
$ws.Cells.Item(67, 8).Value = 1065

# Row 68: Kirguistan
$ws.Cells.Item(68, 2).Value = 48924
$ws.Cells.Item(68, 3).Value = 307
$ws.Cells.Item(68, 4).Value = 44227
$ws.Cells.Item(68, 5).Value = 3615
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 5
$ws.Cells.Item(68, 8).Value = 1082

# Row 69: Ghana
$ws.Cells.Item(69, 2).Value = 47005
$ws.Cells.Item(69, 3).Value = 18
$ws.Cells.Item(69, 4).Value = 46398
$ws.Cells.Item(69, 5).Value = 301
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 306

# Row 112: Haiti
$ws.Cells.Item(112, 2).Value = 8860
$ws.Cells.Item(112, 3).Value = 6
$ws.Cells.Item(112, 4).Value = 7104
$ws.Cells.Item(112, 5).Value = 1526
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 230

# Row 117: Mauritania
$ws.Cells.Item(117, 2).Value = 7550
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(117, 4).Value = 7266
$ws.Cells.Item(117, 5).Value = 121
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 163

# Row 146: Guyana
$ws.Cells.Item(146, 2).Value = 3405
$ws.Cells.Item(146, 3).Value = 47
$ws.Cells.Item(146, 4).Value = 2304
$ws.Cells.Item(146, 5).Value = 999
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 2
$ws.Cells.Item(146, 8).Value = 102

# Row 157: Uruguay
$ws.Cells.Item(157, 2).Value = 2268
$ws.Cells.Item(157, 3).Value = 17
$ws.Cells.Item(157, 4).Value = 1930
$ws.Cells.Item(157, 5).Value = 288
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 1
$ws.Cells.Item(157, 8).Value = 50

# Row 158: Burkina Faso
$ws.Cells.Item(158, 2).Value = 2254
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 1516
$ws.Cells.Item(158, 5).Value = 678
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 60

# Row 191: Barbados
$ws.Cells.Item(191, 2).Value = 206
$ws.Cells.Item(191, 3).Value = 2
$ws.Cells.Item(191, 4).Value = 183
$ws.Cells.Item(191, 5).Value = 16
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 7

